$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Split "Operation Lovecraft: Fallen Doll" into a (playtest) and a new (demo) row ---
# Insert a new row right after the existing "Operation Lovecraft: Fallen Doll" row (209),
# which pushes "Outriders" and everything below down by one row.
$ws.Rows.Item(210).Insert()

# Row 209 keeps the original hash but the title becomes "(playtest)".
$ws.Range("A209").Value = "Operation Lovecraft: Fallen Doll (playtest) 0x64207C05285D224C34D110CB6D935862BB019CC2FE87169E189A97E27A927FAC"

# The newly inserted row 210 is the new "(demo)" entry.
$ws.Range("A210").Value = "Operation Lovecraft: Fallen Doll (demo) 0x496E7699BC1D0DC35DB948772660CF4079F5051408E5DD79E398CF327E6AD6F5"
$ws.Range("B210").Formula = "=LEFT(A210,SEARCH(""/"",SUBSTITUTE(A210,"" "",""/"",LEN(A210)-LEN(SUBSTITUTE(A210,"" "",))))-1)"
$ws.Range("C210").Formula = "=RIGHT(A210,LEN(A210)-FIND(""^^"",SUBSTITUTE(A210,"" "",""^^"",LEN(A210)-LEN(SUBSTITUTE(A210,"" "","""")))))"

# --- Append a new "Fortnite 1.9.1" row at the very end of the table ---
$lastRow = $ws.Cells.Item(1048576, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "Fortnite 1.9.1  DAE1418B289573D4148C72F3C76ABC7E2DB9CAA618A3EAF2D8580EB3A1BB7A63"
$ws.Cells.Item($newRow, 2).Formula = "=LEFT(A" + $newRow + ",SEARCH(""/"",SUBSTITUTE(A" + $newRow + ","" "",""/"",LEN(A" + $newRow + ")-LEN(SUBSTITUTE(A" + $newRow + ","" "",))))-1)"
$ws.Cells.Item($newRow, 3).Formula = "=RIGHT(A" + $newRow + ",LEN(A" + $newRow + ")-FIND(""^^"",SUBSTITUTE(A" + $newRow + ","" "",""^^"",LEN(A" + $newRow + ")-LEN(SUBSTITUTE(A" + $newRow + ","" "","""")))))"

# --- Cosmetic view updates observed in the diff ---
$ws.Range("A1:C" + $newRow).Select()
$excel.ActiveWindow.Zoom = 85
$ws.Application.ActiveWindow.ScrollRow = 281
$ws.Range("C332").Select()
